# Updated cryptos list (GitHub Actions scrape refresh): refresh Price (D)
# and Volume(1h) (E) columns for the latest pull.
#
# Note: several Price values are plain decimals (e.g. "1.00", "534.38").
# Excel's Range.Value setter auto-converts numeric-looking strings to real
# numbers, which would both change the stored cell type and silently drop
# meaningful trailing zeros (e.g. "1.00" -> 1). Prefixing those with a
# leading apostrophe forces Excel to keep them as literal text (matching
# how the source data is stored) without altering the visible text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.563.92"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "2.614.40"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'534.38"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "'142.83"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.576"
$ws.Range("E8").Value = "  +1.74%  "
$ws.Range("E9").Value = "  +3.97%  "
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "'0.136"
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("D13").Value = "3.077.38"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "58.510.78"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "'20.74"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "2.596.12"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "'4.42"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "'334.51"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'66.51"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("D29").Value = "0.0₃0732"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").Value = "'153.44"
$ws.Range("D33").Value = "'18.86"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("D38").Value = "'0.812"
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").Value = "'282.08"
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").Value = "'19.00"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("D46").Value = "'0.0526"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").Value = "1.940.10"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("D50").Value = "'17.83"
$ws.Range("E50").Value = "  -3.61%  "
$ws.Range("D51").Value = "'113.42"
$ws.Range("E51").Value = "  +0.39%  "
